$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1351.3334
$ws.Range("I52").Value = 1402
$ws.Range("J52").Value = 1250
$ws.Range("K52").Value = 4206
$ws.Range("L52").Value = 3750
$ws.Range("M52").Value = -4046
$ws.Range("N52").Value = -4070

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 7753030
$ws.Range("I92").Value = 9804571
$ws.Range("J92").Value = 2764.3333
$ws.Range("K92").Value = 9804571
$ws.Range("L92").Value = 2764.3333
$ws.Range("M92").Value = -9803323
$ws.Range("N92").Value = -5260.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 143591.86
$ws.Range("I103").Value = 500349.5
$ws.Range("J103").Value = 888.8
$ws.Range("K103").Value = 1501048.5
$ws.Range("L103").Value = 2666.4
$ws.Range("M103").Value = -1500462.5
$ws.Range("N103").Value = -3838.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1317.1428
$ws.Range("I2").Value = 1317.1428
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1317.1428
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1204.1428
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 17352.908
$ws.Range("I28").Value = 6485.375
$ws.Range("J28").Value = 46333
$ws.Range("K28").Value = 6485.375
$ws.Range("L28").Value = 46333
$ws.Range("M28").Value = -6293.375
$ws.Range("N28").Value = -46717

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 17898.6
$ws.Range("I43").Value = 8739
$ws.Range("J43").Value = 20188.5
$ws.Range("K43").Value = 8739
$ws.Range("L43").Value = 20188.5
$ws.Range("M43").Value = -8426
$ws.Range("N43").Value = -20814.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 335000
$ws.Range("I97").Value = 500500
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 500500
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -500004
$ws.Range("N97").Value = -4992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 17352.908
$ws.Range("I99").Value = 6485.375
$ws.Range("J99").Value = 46333
$ws.Range("K99").Value = 6485.375
$ws.Range("L99").Value = 46333
$ws.Range("M99").Value = -3490.375
$ws.Range("N99").Value = -52323

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1317.1428
$ws.Range("I116").Value = 1317.1428
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1317.1428
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 976.8571999999999
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1317.1428
$ws.Range("I3").Value = 1317.1428
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1317.1428
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1203.1428
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 599
$ws.Range("I29").Value = 599
$ws.Range("K29").Value = 599
$ws.Range("M29").Value = -310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 4000
$ws.Range("L29").Value = 4000
$ws.Range("N29").Value = -4586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4768.9253
$ws.Range("I31").Value = 1259.625
$ws.Range("K31").Value = 1259.625
$ws.Range("M31").Value = -964.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4768.9253
$ws.Range("I34").Value = 1259.625
$ws.Range("K34").Value = 1259.625
$ws.Range("M34").Value = -1057.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3100.4285
$ws.Range("I86").Value = 3300.6
$ws.Range("K86").Value = 3300.6
$ws.Range("M86").Value = -2177.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3100.4285
$ws.Range("I89").Value = 3300.6
$ws.Range("K89").Value = 16503
$ws.Range("M89").Value = -10887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2126.0857
$ws.Range("I99").Value = 1866.6666
$ws.Range("J99").Value = 2179.7585
$ws.Range("K99").Value = 1866.6666
$ws.Range("L99").Value = 2179.7585
$ws.Range("M99").Value = -368.6666
$ws.Range("N99").Value = -5175.7585

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2126.0857
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 2179.7585
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 6539.2755
$ws.Range("M126").Value = -3129.9998
$ws.Range("N126").Value = -11479.2755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 525.1
$ws.Range("I113").Value = 617.875
$ws.Range("J113").Value = 491.36365
$ws.Range("K113").Value = 1853.625
$ws.Range("L113").Value = 1474.09095
$ws.Range("M113").Value = 316.375
$ws.Range("N113").Value = -5814.09095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4771.276
$ws.Range("I131").Value = 382.14285
$ws.Range("J131").Value = 8867.799999999999
$ws.Range("K131").Value = 1146.42855
$ws.Range("L131").Value = 26603.4
$ws.Range("M131").Value = 3893.57145
$ws.Range("N131").Value = -36683.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 80009
$ws.Range("J22").Value = 80009
$ws.Range("L22").Value = 80009
$ws.Range("N22").Value = -81067

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 64903.562
$ws.Range("I97").Value = 85544.75
$ws.Range("J97").Value = 2980
$ws.Range("K97").Value = 85544.75
$ws.Range("L97").Value = 2980
$ws.Range("M97").Value = -85048.75
$ws.Range("N97").Value = -3972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 735
$ws.Range("J46").Value = 594.44446
$ws.Range("L46").Value = 594.44446
$ws.Range("N46").Value = -970.44446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2950.0571
$ws.Range("I132").Value = 2701.7307
$ws.Range("J132").Value = 3667.4443
$ws.Range("K132").Value = 8105.1921
$ws.Range("L132").Value = 11002.3329
$ws.Range("M132").Value = -5575.1921
$ws.Range("N132").Value = -16062.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1507
$ws.Range("I126").Value = 1698
$ws.Range("J126").Value = 1125
$ws.Range("K126").Value = 5094
$ws.Range("L126").Value = 3375
$ws.Range("M126").Value = -2624
$ws.Range("N126").Value = -8315

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4945574.5
$ws.Range("I132").Value = 1847.0233
$ws.Range("J132").Value = 18231842
$ws.Range("K132").Value = 5541.0699
$ws.Range("L132").Value = 54695526
$ws.Range("M132").Value = -3011.0699
$ws.Range("N132").Value = -54700586
